$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended 2025-05-01T15:37:39.981Z
# Columns: A=ملاحظات B=المرافق C=الكمية D=المخيم E=نوع المسافة F=المركبة G=المؤسسة H=الوقت
$rows = @(
    @{A = $null;  C = "2"; B = "أحمد شريم"; D = "الجزائري"; E = "الرحلة 2"; F = "C2"; G = "NRC"; H = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٧:٣٩ م"},
    @{A = "2";    C = "2"; B = "أحمد شريم"; D = "الجزائري"; E = "الرحلة 2"; F = "C2"; G = "NRC"; H = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٧:٣٩ م"},
    @{A = "2";    C = "2"; B = "أحمد شريم"; D = "الجزائري"; E = "الرحلة 2"; F = "C2"; G = "NRC"; H = "٠١‏/٠٥‏/٢٠٢٥ ٠٦:٣٧:٣٩ م"}
)

$startRow = 6
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]

    # Text columns (never look like numbers) - plain assignment keeps default style.
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
    $ws.Cells.Item($r, 8).Value = $rowData.H

    # Column A: either blank ("" -> no stored value, matches closest achievable state)
    # or a digit-looking value that must be forced to Text so it isn't coerced to a number.
    if ($rowData.A) {
        $cellA = $ws.Cells.Item($r, 1)
        $cellA.NumberFormat = "@"
        $cellA.Value = $rowData.A
    }

    # Column C ("الكمية") always holds a digit-looking value here - force Text so the
    # stored cell stays a string ("2"), matching the numberStoredAsText source data.
    $cellC = $ws.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $rowData.C
}
